$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second thumbnail entry (row 3, "thumbnail2.jpg") as part of
# implementing the add-images logic against the excel data.
$ws.Range("C3").ClearContents()

# Leave the selection on the cell that was just edited.
$ws.Range("C3").Select() | Out-Null
